# Resolve War Packing issues: remove the extra flight/PNR columns (D:H)
# that had been duplicated into rows 3-17 of the generated test case sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3:H17").ClearContents()
